$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Segundo cambio en text branch"
$ws.Range("A2").Select()
